# phone type clean up
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scouting Admin")

# --- Field Questions section ---
$ws.Range("A12").Copy($ws.Range("A25"))
$ws.Range("B12").Copy($ws.Range("B25"))
$ws.Range("C12").Copy($ws.Range("C25"))
$ws.Range("A25").Value = "Field Questions"
$ws.Range("A25:C25").Merge()

$ws.Range("A13").Copy($ws.Range("A26"))
$ws.Range("B13").Copy($ws.Range("B26"))
$ws.Range("A26").Value = "Save"

$ws.Range("A13").Copy($ws.Range("A27"))
$ws.Range("B13").Copy($ws.Range("B27"))
$ws.Range("A27").Value = "Edit"

# --- Pit Questions section ---
$ws.Range("A12").Copy($ws.Range("A28"))
$ws.Range("B12").Copy($ws.Range("B28"))
$ws.Range("C12").Copy($ws.Range("C28"))
$ws.Range("A28").Value = "Pit Questions"
$ws.Range("A28:C28").Merge()

$ws.Range("A13").Copy($ws.Range("A29"))
$ws.Range("B13").Copy($ws.Range("B29"))
$ws.Range("A29").Value = "Save"

$ws.Range("A13").Copy($ws.Range("A30"))
$ws.Range("B13").Copy($ws.Range("B30"))
$ws.Range("A30").Value = "Edit"

# --- Phone Types section ---
$ws.Range("A12").Copy($ws.Range("A31"))
$ws.Range("B12").Copy($ws.Range("B31"))
$ws.Range("C12").Copy($ws.Range("C31"))
$ws.Range("A31").Value = "Phone Types"
$ws.Range("A31:C31").Merge()

$ws.Range("A13").Copy($ws.Range("A32"))
$ws.Range("B13").Copy($ws.Range("B32"))
$ws.Range("A32").Value = "Save"

$ws.Range("A13").Copy($ws.Range("A33"))
$ws.Range("B13").Copy($ws.Range("B33"))
$ws.Range("A33").Value = "Edit"

$ws.Range("A32:B33").Select()
